$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1183.3
$ws.Range("J112").Value = 1183.3
$ws.Range("L112").Value = 3549.9
$ws.Range("N112").Value = -5765.9

$ws.Range("H116").Value = 6673566.5
$ws.Range("I116").Value = 9850
$ws.Range("J116").Value = 11116044
$ws.Range("K116").Value = 9850
$ws.Range("L116").Value = 11116044
$ws.Range("M116").Value = -6408
$ws.Range("N116").Value = -11122928

$ws.Range("H132").Value = 1639.409
$ws.Range("I132").Value = 1462.5883
$ws.Range("K132").Value = 4387.7649
$ws.Range("M132").Value = -1857.7649

$ws.Range("H137").Value = 485907.34
$ws.Range("I137").Value = 1875.1578
$ws.Range("K137").Value = 5625.4734
$ws.Range("M137").Value = -3075.4734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 169622.33
$ws.Range("I61").Value = 3547
$ws.Range("K61").Value = 3547
$ws.Range("M61").Value = -3335

$ws.Range("H74").Value = 44333.207
$ws.Range("I74").Value = 68347
$ws.Range("K74").Value = 68347
$ws.Range("M74").Value = -67473

$ws.Range("H77").Value = 44333.207
$ws.Range("I77").Value = 68347
$ws.Range("K77").Value = 341735
$ws.Range("M77").Value = -337367

$ws.Range("H132").Value = 2326.5151
$ws.Range("I132").Value = 2102.963
$ws.Range("J132").Value = 3332.5
$ws.Range("K132").Value = 6308.889000000001
$ws.Range("L132").Value = 9997.5
$ws.Range("M132").Value = -3778.889000000001
$ws.Range("N132").Value = -15057.5

$ws.Range("H136").Value = 169622.33
$ws.Range("I136").Value = 3547
$ws.Range("K136").Value = 10641
$ws.Range("M136").Value = -8091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1359701.6
$ws.Range("I7").Value = 1811275.8
$ws.Range("J7").Value = 4979.3335
$ws.Range("K7").Value = 1811275.8
$ws.Range("L7").Value = 4979.3335
$ws.Range("M7").Value = -1811162.8
$ws.Range("N7").Value = -5205.3335

$ws.Range("H134").Value = 5249.722
$ws.Range("I134").Value = 2891.5715
$ws.Range("J134").Value = 13503.25
$ws.Range("K134").Value = 8674.7145
$ws.Range("L134").Value = 40509.75
$ws.Range("M134").Value = -6139.7145
$ws.Range("N134").Value = -45579.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 152.14285
$ws.Range("J2").Value = 101.25
$ws.Range("L2").Value = 101.25
$ws.Range("N2").Value = -327.25

$ws.Range("H31").Value = 2780.258
$ws.Range("I31").Value = 1770
$ws.Range("K31").Value = 1770
$ws.Range("M31").Value = -1475

$ws.Range("H34").Value = 2780.258
$ws.Range("I34").Value = 1770
$ws.Range("K34").Value = 1770
$ws.Range("M34").Value = -1568

$ws.Range("H94").Value = 898.2
$ws.Range("I94").Value = 899.6667
$ws.Range("J94").Value = 896
$ws.Range("K94").Value = 899.6667
$ws.Range("L94").Value = 896
$ws.Range("M94").Value = -448.6667
$ws.Range("N94").Value = -1798

$ws.Range("H132").Value = 449817.8
$ws.Range("I132").Value = 1609.4073
$ws.Range("K132").Value = 4828.2219
$ws.Range("M132").Value = -2298.2219

$ws.Range("H134").Value = 65570
$ws.Range("I134").Value = 3305.5833
$ws.Range("J134").Value = 252363.25
$ws.Range("K134").Value = 9916.749899999999
$ws.Range("L134").Value = 757089.75
$ws.Range("M134").Value = -7381.749899999999
$ws.Range("N134").Value = -762159.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3136.2
$ws.Range("I3").Value = 3136.2
$ws.Range("K3").Value = 9408.599999999999
$ws.Range("M3").Value = -9296.599999999999

$ws.Range("H8").Value = 130133.4
$ws.Range("I8").Value = 130133.4
$ws.Range("K8").Value = 390400.2
$ws.Range("M8").Value = -390261.2

$ws.Range("H13").Value = 500
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""

$ws.Range("H36").Value = 1076.25
$ws.Range("I36").Value = 1094.2858
$ws.Range("J36").Value = 950
$ws.Range("K36").Value = 3282.8574
$ws.Range("L36").Value = 2850
$ws.Range("M36").Value = -3113.8574
$ws.Range("N36").Value = -3188

$ws.Range("H68").Value = 102169
$ws.Range("J68").Value = 113243.336
$ws.Range("L68").Value = 339730.008
$ws.Range("N68").Value = -341352.008

$ws.Range("H71").Value = 102169
$ws.Range("J71").Value = 113243.336
$ws.Range("L71").Value = 1019190.024
$ws.Range("N71").Value = -1027302.024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1258.8667
$ws.Range("I102").Value = 1242
$ws.Range("K102").Value = 1242
$ws.Range("M102").Value = 380

$ws.Range("H122").Value = 8634989
$ws.Range("I122").Value = 9354482
$ws.Range("J122").Value = 1075
$ws.Range("K122").Value = 28063446
$ws.Range("L122").Value = 3225
$ws.Range("M122").Value = -28060996
$ws.Range("N122").Value = -8125

$ws.Range("H132").Value = 3665.5454
$ws.Range("I132").Value = 2546.4285
$ws.Range("K132").Value = 7639.2855
$ws.Range("M132").Value = -5109.2855

$ws.Range("H135").Value = 39899.637
$ws.Range("J135").Value = 39899.637
$ws.Range("L135").Value = 39899.637
$ws.Range("N135").Value = -50039.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1889.6
$ws.Range("I16").Value = 2873
$ws.Range("J16").Value = 906.2
$ws.Range("K16").Value = 2873
$ws.Range("L16").Value = 906.2
$ws.Range("M16").Value = -2703
$ws.Range("N16").Value = -1246.2

$ws.Range("H46").Value = 6046.136
$ws.Range("I46").Value = 7735.375
$ws.Range("J46").Value = 1541.5
$ws.Range("K46").Value = 7735.375
$ws.Range("L46").Value = 1541.5
$ws.Range("M46").Value = -7547.375
$ws.Range("N46").Value = -1917.5

$ws.Range("H136").Value = 3086.724
$ws.Range("I136").Value = 3209.1667
$ws.Range("J136").Value = 2499
$ws.Range("K136").Value = 9627.500100000001
$ws.Range("L136").Value = 7497
$ws.Range("M136").Value = -7077.500100000001
$ws.Range("N136").Value = -12597

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 98304.586
$ws.Range("I126").Value = 162808
$ws.Range("J126").Value = 7999.8
$ws.Range("K126").Value = 488424
$ws.Range("L126").Value = 23999.4
$ws.Range("M126").Value = -485954
$ws.Range("N126").Value = -28939.4

$ws.Range("H132").Value = 1704.8628
$ws.Range("I132").Value = 1429.5135
$ws.Range("J132").Value = 2432.5715
$ws.Range("K132").Value = 4288.5405
$ws.Range("L132").Value = 7297.7145
$ws.Range("M132").Value = -1758.5405
